# CU Consultar proximos pagos: ligar un pago a un grupo.
# "PagoAlumno" relation gains an "idGrupo" attribute just before the
# closing brace, e.g.
#   PagoAlumno {idPago, idAlumno, idPromocion, tipoPago, monto, fecha}.
# becomes
#   PagoAlumno {idPago, idAlumno, idPromocion, tipoPago, monto, fecha, idGrupo}.

$d = $word.ActiveDocument

# Scope the Find/Replace to the first paragraph only, since the sequence
# ", monto, fecha" also appears (differently punctuated) later in the
# document (e.g. the "Egreso" relation) and must stay untouched.
$rng = $d.Paragraphs(1).Range

$found = $rng.Find.Execute(", monto, fecha}.", $true, $false, $false, $false, $false, $true, 1, $false, ", monto, fecha, idGrupo}.", 2)

if (-not $found) {
    Write-Output "WARNING: target text not found in first paragraph"
}

Write-Output $d.Paragraphs(1).Range.Text
